$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 updates
$ws.Range("G4").Value = 3.5
$ws.Range("H4").Value = 2.8
$ws.Range("I4").Value = 2.98
$ws.Range("J4").Value = 2.82
$ws.Range("P4").Value = 1.3
$ws.Range("T4").Value = 2.68
$ws.Range("Y4").Value = 1000
$ws.Range("Z4").Value = 20
$ws.Range("AJ4").Value = 1000

# Row 7 updates
$ws.Range("N7").Value = 2.46
$ws.Range("P7").Value = 1.48
$ws.Range("Q7").Value = 2.74
$ws.Range("AE7").Value = 85
$ws.Range("AN7").Value = 980
